$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 857.2857
$ws.Range("I28").Value = 880.3333
$ws.Range("J28").Value = 840
$ws.Range("K28").Value = 880.3333
$ws.Range("L28").Value = 840
$ws.Range("M28").Value = -395.3333
$ws.Range("N28").Value = -1810

$ws.Range("H61").Value = 80
$ws.Range("I61").Value = 80
$ws.Range("K61").Value = 240
$ws.Range("M61").Value = -68

$ws.Range("H132").Value = 2437.8572
$ws.Range("I132").Value = 1073.1
$ws.Range("K132").Value = 3219.3
$ws.Range("M132").Value = -689.2999999999997

$ws.Range("H135").Value = 492.5
$ws.Range("I135").Value = 492.5
$ws.Range("K135").Value = 4432.5
$ws.Range("M135").Value = -1897.5

$ws.Range("H137").Value = 3440.25
$ws.Range("I137").Value = 2966.3333
$ws.Range("J137").Value = 3598.2222
$ws.Range("K137").Value = 8898.999899999999
$ws.Range("L137").Value = 10794.6666
$ws.Range("M137").Value = -6348.999899999999
$ws.Range("N137").Value = -15894.6666

$ws.Range("H138").Value = 13330.333
$ws.Range("I138").Value = 999.3333
$ws.Range("K138").Value = 2997.9999
$ws.Range("M138").Value = 2142.0001

$ws.Range("H141").Value = 4499.5
$ws.Range("I141").Value = 4499.5
$ws.Range("K141").Value = 13498.5
$ws.Range("M141").Value = -8318.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1649.9
$ws.Range("I4").Value = 1800
$ws.Range("J4").Value = 299
$ws.Range("K4").Value = 1800
$ws.Range("L4").Value = 299
$ws.Range("M4").Value = -1684
$ws.Range("N4").Value = -531

$ws.Range("H32").Value = 9645
$ws.Range("I32").Value = 8859.412
$ws.Range("K32").Value = 8859.412
$ws.Range("M32").Value = -8572.412

$ws.Range("M74").ClearContents()
$ws.Range("H74").Value = 20000
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = -19126

$ws.Range("M77").ClearContents()
$ws.Range("H77").Value = 20000
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 100000
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = -95632

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 491.25
$ws.Range("I22").Value = 490
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 490
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -317
$ws.Range("N22").Value = -846

$ws.Range("H86").Value = 11738.667
$ws.Range("I86").Value = 4787.5
$ws.Range("K86").Value = 4787.5
$ws.Range("M86").Value = -3664.5

$ws.Range("H89").Value = 11738.667
$ws.Range("I89").Value = 4787.5
$ws.Range("K89").Value = 23937.5
$ws.Range("M89").Value = -18321.5

$ws.Range("H99").Value = 2597
$ws.Range("I99").Value = 2421.625
$ws.Range("K99").Value = 2421.625
$ws.Range("M99").Value = -923.625

$ws.Range("H105").Value = 7874.6665
$ws.Range("I105").Value = 11142.857
$ws.Range("J105").Value = 3299.2
$ws.Range("K105").Value = 11142.857
$ws.Range("L105").Value = 3299.2
$ws.Range("M105").Value = -9395.857
$ws.Range("N105").Value = -6793.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 411.64706
$ws.Range("I22").Value = 399.92856
$ws.Range("J22").Value = 466.33334
$ws.Range("K22").Value = 399.92856
$ws.Range("L22").Value = 466.33334
$ws.Range("M22").Value = -49.92856
$ws.Range("N22").Value = -1166.33334

$ws.Range("H105").Value = 2632.8333
$ws.Range("I105").Value = 2632.8333
$ws.Range("K105").Value = 2632.8333
$ws.Range("M105").Value = -885.8332999999998

$ws.Range("H132").Value = 147642.42
$ws.Range("J132").Value = 6374.5
$ws.Range("L132").Value = 19123.5
$ws.Range("N132").Value = -24183.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L68").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("N68").Value = 0

$ws.Range("L71").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("N71").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4749.5
$ws.Range("I132").Value = 3999.5
$ws.Range("K132").Value = 11998.5
$ws.Range("M132").Value = -9468.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10205.077
$ws.Range("J22").Value = 8583.333000000001
$ws.Range("L22").Value = 8583.333000000001
$ws.Range("N22").Value = -9173.333000000001

$ws.Range("H27").Value = 10205.077
$ws.Range("J27").Value = 8583.333000000001
$ws.Range("L27").Value = 8583.333000000001
$ws.Range("N27").Value = -8797.333000000001

$ws.Range("H46").Value = 3312.5
$ws.Range("J46").Value = 3125
$ws.Range("L46").Value = 3125
$ws.Range("N46").Value = -3501

$ws.Range("H55").Value = 561.8
$ws.Range("I55").Value = 636.4167
$ws.Range("J55").Value = 263.33334
$ws.Range("K55").Value = 636.4167
$ws.Range("L55").Value = 263.33334
$ws.Range("M55").Value = -463.4167
$ws.Range("N55").Value = -609.33334

$ws.Range("H61").Value = 12752149
$ws.Range("I61").Value = 10202418
$ws.Range("K61").Value = 10202418
$ws.Range("M61").Value = -10202216

$ws.Range("H93").Value = 2250
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -752
$ws.Range("N93").Value = -4996

$ws.Range("H113").Value = 12752149
$ws.Range("I113").Value = 10202418
$ws.Range("K113").Value = 10202418
$ws.Range("M113").Value = -10200248

$ws.Range("H132").Value = 5199.4
$ws.Range("I132").Value = 3999.5
$ws.Range("K132").Value = 11998.5
$ws.Range("M132").Value = -9468.5

$ws.Range("H136").Value = 5930.4546
$ws.Range("I136").Value = 2716.875
$ws.Range("K136").Value = 8150.625
$ws.Range("M136").Value = -5600.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6301176
$ws.Range("I2").Value = 6382499.5
$ws.Range("K2").Value = 6382499.5
$ws.Range("M2").Value = -6382387.5

$ws.Range("H8").Value = 8999.5
$ws.Range("J8").Value = 8999.5
$ws.Range("L8").Value = 8999.5
$ws.Range("N8").Value = -9279.5

$ws.Range("H107").Value = 623.5714
$ws.Range("I107").Value = 668.875
$ws.Range("J107").Value = 563.1667
$ws.Range("K107").Value = 2006.625
$ws.Range("L107").Value = 1689.5001
$ws.Range("M107").Value = -86.625
$ws.Range("N107").Value = -5529.5001

$ws.Range("H113").Value = 706.17645
$ws.Range("I113").Value = 621.8333
$ws.Range("K113").Value = 1865.4999
$ws.Range("M113").Value = 304.5001

$ws.Range("H122").Value = 2135.4546
$ws.Range("I122").Value = 2898.5715
$ws.Range("K122").Value = 8695.7145
$ws.Range("M122").Value = -6245.7145

$ws.Range("H132").Value = 4349.5
$ws.Range("I132").Value = 3999.6667
$ws.Range("K132").Value = 11999.0001
$ws.Range("M132").Value = -9469.000100000001
